# Auto-generated Excel COM-interop script applying the crypto price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.930.28'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '1.632.52'
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.51'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.22%  '
$ws.Range('E9').Value = '  +1.41%  '
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').Value = '1.866.27'
$ws.Range('E12').Value = '  +1.77%  '
$ws.Range('D13').Value = '1.626.30'
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('E14').Value = '  +2.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +17.74%  '
$ws.Range('D16').Value = '29.941.19'
$ws.Range('E16').Value = '  +0.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.85'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').Value = '0.0₃0701'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.84'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.86%  '
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.15'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.50'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('E28').Value = '  +2.42%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('E30').Value = '  +1.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.11'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.06%  '
$ws.Range('E32').Value = '  +4.12%  '
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('D34').Value = '1.423.43'
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('E35').Value = '  +4.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.03'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.80'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.85%  '
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '75.76'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +12.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.551'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.72%  '
$ws.Range('E42').Value = '  +1.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.826'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0488'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.76%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.02'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '53.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.86%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('D48').Value = '1.774.49'
$ws.Range('E48').Value = '  +1.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.34'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0113'
$ws.Range('E50').Value = '  +9.61%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '89.62'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.76%  '
